# HIKER-M Update Attendance List [TV]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The two recorded absences now use the "heavy ballot X" glyph instead of a
# plain "x" (sharedStrings gains "✘" and drops the old "x" entry).
$ws.Range("H8").Value = "✘"
$ws.Range("H14").Value = "✘"

# A new attendance-taking session (2021-06-02) was added as column L.
$ws.Range("L4").Value = 44349
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats - match the date-header style

$attendance = @{
    5  = "✓ (Joined at 14:30)"   # Feldgrill joined the session late
    6  = "✓"
    7  = "✓"
    8  = "✓"
    9  = "✓"
    10 = "✓"
    11 = "✓"
    12 = "✓"
    13 = "✓"
    14 = "✘"                      # Uka was absent, same as H14
    15 = "✓"
}

foreach ($row in $attendance.Keys | Sort-Object) {
    $target = $ws.Range("L$row")
    $target.Value = $attendance[$row]

    if ($row -eq 14) {
        # Reuse the "absent" formatting already applied to H14 rather than
        # the usual checkmark style used by the rest of column L.
        $ws.Range("H14").Copy()
    }
    else {
        $ws.Range("K$row").Copy()
    }
    $target.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# New column needs its own width (23.5 chars, to fit "✓ (Joined at 14:30)").
$ws.Columns.Item(12).ColumnWidth = 22.666666666666668

$ws.Range("L9").Select()
